$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 118: Crafty Concoctions / Commanding Craftsman's Syrup
$ws.Range("H118").Value = 1322.2222
$ws.Range("I118").Value = 300
$ws.Range("J118").Value = 2600
$ws.Range("K118").Value = 900
$ws.Range("L118").Value = 7800
$ws.Range("M118").Value = 757
$ws.Range("N118").Value = -11114

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 4960.9165
$ws.Range("I2").Value = 4277.5
$ws.Range("J2").Value = 6327.75
$ws.Range("K2").Value = 4277.5
$ws.Range("L2").Value = 6327.75
$ws.Range("M2").Value = -4164.5
$ws.Range("N2").Value = -6553.75
# Row 5: The Alloyed Truth / Bronze Rivets
$ws.Range("H5").Value = 175
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = 62
$ws.Range("N5").Value = -524
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 1459.8572
$ws.Range("I61").Value = 1390.1111
$ws.Range("J61").Value = 1585.4
$ws.Range("K61").Value = 1390.1111
$ws.Range("L61").Value = 1585.4
$ws.Range("M61").Value = -1178.1111
$ws.Range("N61").Value = -2009.4
# Row 97: Ore for Me / High Steel Ingot
$ws.Range("H97").Value = 2000.9166
$ws.Range("I97").Value = 1980
$ws.Range("K97").Value = 1980
$ws.Range("M97").Value = -1484
# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 2222.1304
$ws.Range("I102").Value = 1247.2667
$ws.Range("J102").Value = 4050
$ws.Range("K102").Value = 1247.2667
$ws.Range("L102").Value = 4050
$ws.Range("M102").Value = 374.7333000000001
$ws.Range("N102").Value = -7294
# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 1151.1
$ws.Range("I110").Value = 1151.1
$ws.Range("K110").Value = 1151.1
$ws.Range("M110").Value = 893.9000000000001
# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 4960.9165
$ws.Range("I116").Value = 4277.5
$ws.Range("J116").Value = 6327.75
$ws.Range("K116").Value = 4277.5
$ws.Range("L116").Value = 6327.75
$ws.Range("M116").Value = -1983.5
$ws.Range("N116").Value = -10915.75
# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 34642.773
$ws.Range("I122").Value = 2286.8333
$ws.Range("K122").Value = 6860.499899999999
$ws.Range("M122").Value = -4410.499899999999
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 43524550
$ws.Range("I132").Value = 90911520
$ws.Range("J132").Value = 86497.664
$ws.Range("K132").Value = 272734560
$ws.Range("L132").Value = 259492.992
$ws.Range("M132").Value = -272732030
$ws.Range("N132").Value = -264552.992
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 1459.8572
$ws.Range("I136").Value = 1390.1111
$ws.Range("J136").Value = 1585.4
$ws.Range("K136").Value = 4170.3333
$ws.Range("L136").Value = 4756.200000000001
$ws.Range("M136").Value = -1620.3333
$ws.Range("N136").Value = -9856.200000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 4960.9165
$ws.Range("I3").Value = 4277.5
$ws.Range("J3").Value = 6327.75
$ws.Range("K3").Value = 4277.5
$ws.Range("L3").Value = 6327.75
$ws.Range("M3").Value = -4163.5
$ws.Range("N3").Value = -6555.75
# Row 4: Mending Fences / Bronze Rivets
$ws.Range("H4").Value = 175
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 50
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = 65
$ws.Range("N4").Value = -530
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 12333.333
$ws.Range("I20").Value = 12000
$ws.Range("J20").Value = 12500
$ws.Range("K20").Value = 12000
$ws.Range("L20").Value = 12500
$ws.Range("M20").Value = -11753
$ws.Range("N20").Value = -12994
# Row 24: Honest Ballast / Initiate's Head Knife
$ws.Range("H24").Value = 1371.4
$ws.Range("I24").Value = 464.25
$ws.Range("K24").Value = 464.25
$ws.Range("M24").Value = -229.25
# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1681.1
$ws.Range("I105").Value = 1506.6666
$ws.Range("J105").Value = 1755.8572
$ws.Range("K105").Value = 1506.6666
$ws.Range("L105").Value = 1755.8572
$ws.Range("M105").Value = 240.3334
$ws.Range("N105").Value = -5249.8572

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 5815825
$ws.Range("I31").Value = 7144516.5
$ws.Range("J31").Value = 2800
$ws.Range("K31").Value = 7144516.5
$ws.Range("L31").Value = 2800
$ws.Range("M31").Value = -7144221.5
$ws.Range("N31").Value = -3390
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 5815825
$ws.Range("I34").Value = 7144516.5
$ws.Range("J34").Value = 2800
$ws.Range("K34").Value = 7144516.5
$ws.Range("L34").Value = 2800
$ws.Range("M34").Value = -7144314.5
$ws.Range("N34").Value = -3204
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 2900
$ws.Range("J62").Value = 2900
$ws.Range("L62").Value = 2900
$ws.Range("N62").Value = -4148
# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 2900
$ws.Range("J65").Value = 2900
$ws.Range("L65").Value = 14500
$ws.Range("N65").Value = -20740
# Row 86: Birch, Please / Birch Lumber
$ws.Range("H86").Value = 36925.7
$ws.Range("I86").Value = 1865.2858
$ws.Range("J86").Value = 118733.336
$ws.Range("K86").Value = 1865.2858
$ws.Range("L86").Value = 118733.336
$ws.Range("M86").Value = -742.2858000000001
$ws.Range("N86").Value = -120979.336
# Row 89: Built This City on Blocks and Soul (L) / Birch Lumber
$ws.Range("H89").Value = 36925.7
$ws.Range("I89").Value = 1865.2858
$ws.Range("J89").Value = 118733.336
$ws.Range("K89").Value = 9326.429
$ws.Range("L89").Value = 593666.6799999999
$ws.Range("M89").Value = -3710.429
$ws.Range("N89").Value = -604898.6799999999
# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 1325.25
$ws.Range("I105").Value = 1025.3334
$ws.Range("J105").Value = 2225
$ws.Range("K105").Value = 1025.3334
$ws.Range("L105").Value = 2225
$ws.Range("M105").Value = 721.6666
$ws.Range("N105").Value = -5719
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 62511.65
$ws.Range("I132").Value = 2522.2222
$ws.Range("J132").Value = 129999.75
$ws.Range("K132").Value = 7566.6666
$ws.Range("L132").Value = 389999.25
$ws.Range("M132").Value = -5036.6666
$ws.Range("N132").Value = -395059.25

$ws = $wb.Worksheets.Item("CUL")
# Row 113: Can't Eat Just One / Night Vinegar
$ws.Range("H113").Value = 973.0513
$ws.Range("I113").Value = 572
$ws.Range("J113").Value = 1006.4722
$ws.Range("K113").Value = 1716
$ws.Range("L113").Value = 3019.4166
$ws.Range("M113").Value = 454
$ws.Range("N113").Value = -7359.4166

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 26794500
$ws.Range("I70").Value = 51144330
$ws.Range("J70").Value = 9691
$ws.Range("K70").Value = 51144330
$ws.Range("L70").Value = 9691
$ws.Range("M70").Value = -51144060
$ws.Range("N70").Value = -10231
# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 26794500
$ws.Range("I73").Value = 51144330
$ws.Range("J73").Value = 9691
$ws.Range("K73").Value = 51144330
$ws.Range("L73").Value = 9691
$ws.Range("M73").Value = -51143394
$ws.Range("N73").Value = -11563

$ws = $wb.Worksheets.Item("LTW")
# Row 25: A Rush on Ringbands / Hard Leather Ringbands
$ws.Range("H25").Value = 1366
$ws.Range("I25").Value = 533.3333
$ws.Range("J25").Value = 2198.6667
$ws.Range("K25").Value = 533.3333
$ws.Range("L25").Value = 2198.6667
$ws.Range("M25").Value = -303.3333
$ws.Range("N25").Value = -2658.6667
# Row 133: The Perfect Accessory / Loboskin Amulet of Fending
$ws.Range("H133").Value = 113400
$ws.Range("J133").Value = 113400
$ws.Range("L133").Value = 113400
$ws.Range("N133").Value = -118460

$ws = $wb.Worksheets.Item("WVR")
# Row 111: Legs for Days / Iridescent Bottoms of Scouting
$ws.Range("H111").Value = 24999.95
$ws.Range("J111").Value = 24999.95
$ws.Range("L111").Value = 24999.95
$ws.Range("N111").Value = -33179.95
# Row 133: Begin with the Basics / Snow Cotton Jacket
$ws.Range("H133").Value = 42857.5
$ws.Range("J133").Value = 42857.5
$ws.Range("L133").Value = 42857.5
$ws.Range("N133").Value = -52977.5
